$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row (row 17) below the existing worker record (row 16),
# copying that row's full formatting so the new row matches the table style.
$ws.Rows(16).Copy()
$ws.Rows(17).Insert()
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# The new row represents the same worker for the next "Periodo Mora" (2509).
$ws.Range("E17").Value = "2509"

# Update the summary counters/totals to reflect the second period being added.
$ws.Range("F13").Value = 2
$ws.Range("E11").Value = 113880

$excel.CutCopyMode = 0
